$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$rng = $ws.Range("Z2")
$rng.Value = "Test"
$rng.Borders.Item(7).LineStyle = 1
$rng.Borders.Item(7).ColorIndex = 1
